$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Project Charge Code" value for the sample row was renamed from
# "Pro-Lot Track (Lot Track)" to "Pro-SYDATA1 (Lot track)".
$ws.Range("B2").Value = "Pro-SYDATA1 (Lot track)"

# Update the saved selection/active cell on the sheet to just B2.
$ws.Range("B2").Select()

$wb.Save()
